$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.906.98'
$ws.Range('E2').Value = '  +2.29%  '
$ws.Range('D3').Value = '1.810.06'
$ws.Range('E3').Value = '  +2.97%  '
$ws.Range('E4').Value = '  +0.49%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.02'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.006'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.48%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4290'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3686'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07240'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +3.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8634'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.33%  '
$ws.Range('D11').Value = '2.041.33'
$ws.Range('E11').Value = '  +17.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.19'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +5.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.621'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +4.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.393'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.68%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06930'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '80.81'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.012'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008845'
$ws.Range('D18').ClearFormats()
$ws.Range('E19').Value = '  +0.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.17'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.34%  '
$ws.Range('D21').Value = '26.937.58'
$ws.Range('E21').Value = '  +2.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.192'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +4.61%  '
$ws.Range('E23').Value = '  -0.58%  '
$ws.Range('D24').Value = '2.279.09'
$ws.Range('E24').Value = '  +16.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.02'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.883'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.31'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.240'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.91%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.912'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +14.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '114.65'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08948'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7396'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.158'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.428'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.806'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.76%  '
$ws.Range('E36').Value = '  +0.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.123'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +5.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05222'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01923'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5082'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +4.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.761'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +11.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1647'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.455'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +5.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.293'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +5.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '107.21'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.70%  '
$ws.Range('E46').Value = '  +4.14%  '
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.644'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +5.39%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4551'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +2.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06271'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.89%  '
$ws.Range('E51').Value = '  +6.30%  '
